$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Clean up main/data input: update clr (L2) and cnb (N2) values
$ws.Range("L2").Value = 0.01
$ws.Range("N2").Value = 0.085
